# ECS-62: update report header labels for the experiment report template
# "Дата начала эксперимента" / "Дата окончания эксперимента" ->
# "Дата начала обработки заявки" / "Дата окончания обработки заявки"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Отчет по заявкам на эксперимент")

$ws.Range("K8").Value = "Дата начала обработки заявки"
$ws.Range("L8").Value = "Дата окончания обработки заявки"
